$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 117.42857
$ws.Range("I33").Value = 122.31579
$ws.Range("K33").Value = 122.31579
$ws.Range("M33").Value = 106.68421
$ws.Range("H93").Value = 29601
$ws.Range("J93").Value = 29601
$ws.Range("L93").Value = 29601
$ws.Range("N93").Value = -34593
$ws.Range("H137").Value = 2777.682
$ws.Range("I137").Value = 2694.647
$ws.Range("J137").Value = 3060
$ws.Range("K137").Value = 8083.941
$ws.Range("L137").Value = 9180
$ws.Range("M137").Value = -5533.941
$ws.Range("N137").Value = -14280
$ws.Range("H138").Value = 1781.7625
$ws.Range("I138").Value = 718.8889
$ws.Range("J138").Value = 2323.2263
$ws.Range("K138").Value = 2156.6667
$ws.Range("L138").Value = 6969.678899999999
$ws.Range("M138").Value = 2983.3333
$ws.Range("N138").Value = -17249.6789
$ws.Range("H141").Value = 1980
$ws.Range("I141").Value = 1693.3334
$ws.Range("J141").Value = 2625
$ws.Range("K141").Value = 5080.0002
$ws.Range("L141").Value = 7875
$ws.Range("M141").Value = 99.9997999999996
$ws.Range("N141").Value = -18235
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1245.9412
$ws.Range("I2").Value = 1210.2273
$ws.Range("K2").Value = 1210.2273
$ws.Range("M2").Value = -1097.2273
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -885
$ws.Range("H63").Value = 7814124
$ws.Range("J63").Value = 2500
$ws.Range("L63").Value = 2500
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 7814124
$ws.Range("J66").Value = 2500
$ws.Range("L66").Value = 12500
$ws.Range("N66").Value = -19364
$ws.Range("H92").Value = 23000
$ws.Range("J92").Value = 23000
$ws.Range("L92").Value = 23000
$ws.Range("N92").Value = -27992
$ws.Range("H112").Value = 36999.5
$ws.Range("J112").Value = 36999.5
$ws.Range("L112").Value = 36999.5
$ws.Range("N112").Value = -39953.5
$ws.Range("H114").Value = 27016.5
$ws.Range("J114").Value = 27016.5
$ws.Range("L114").Value = 27016.5
$ws.Range("N114").Value = -35694.5
$ws.Range("H116").Value = 1245.9412
$ws.Range("I116").Value = 1210.2273
$ws.Range("K116").Value = 1210.2273
$ws.Range("M116").Value = 1083.7727
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1245.9412
$ws.Range("I3").Value = 1210.2273
$ws.Range("K3").Value = 1210.2273
$ws.Range("M3").Value = -1096.2273
$ws.Range("H92").Value = 29899.334
$ws.Range("J92").Value = 29899.334
$ws.Range("L92").Value = 29899.334
$ws.Range("N92").Value = -34891.334
$ws.Range("H110").Value = 45189.2
$ws.Range("J110").Value = 45189.2
$ws.Range("L110").Value = 45189.2
$ws.Range("N110").Value = -53369.2
$ws.Range("H134").Value = 3085.7632
$ws.Range("I134").Value = 3316.4412
$ws.Range("J134").Value = 1125
$ws.Range("K134").Value = 9949.3236
$ws.Range("L134").Value = 3375
$ws.Range("M134").Value = -7414.3236
$ws.Range("N134").Value = -8445
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10895.6875
$ws.Range("I31").Value = 15020.655
$ws.Range("K31").Value = 15020.655
$ws.Range("M31").Value = -14725.655
$ws.Range("H34").Value = 10895.6875
$ws.Range("I34").Value = 15020.655
$ws.Range("K34").Value = 15020.655
$ws.Range("M34").Value = -14818.655
$ws.Range("H124").Value = 37733.332
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H132").Value = 12187.16
$ws.Range("I132").Value = 16346.735
$ws.Range("J132").Value = 3348.0625
$ws.Range("K132").Value = 49040.205
$ws.Range("L132").Value = 10044.1875
$ws.Range("M132").Value = -46510.205
$ws.Range("N132").Value = -15104.1875
$ws.Range("H134").Value = 1108.3871
$ws.Range("I134").Value = 846.3077
$ws.Range("J134").Value = 1552.7826
$ws.Range("K134").Value = 2538.9231
$ws.Range("L134").Value = 4658.3478
$ws.Range("M134").Value = -3.923099999999977
$ws.Range("N134").Value = -9728.3478
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 773.55
$ws.Range("J131").Value = 773.55
$ws.Range("L131").Value = 2320.65
$ws.Range("N131").Value = -12400.65
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3091.4644
$ws.Range("I80").Value = 2664.5833
$ws.Range("J80").Value = 3411.625
$ws.Range("K80").Value = 2664.5833
$ws.Range("L80").Value = 3411.625
$ws.Range("M80").Value = -1666.5833
$ws.Range("N80").Value = -5407.625
$ws.Range("H83").Value = 3091.4644
$ws.Range("I83").Value = 2664.5833
$ws.Range("J83").Value = 3411.625
$ws.Range("K83").Value = 13322.9165
$ws.Range("L83").Value = 17058.125
$ws.Range("M83").Value = -8330.916499999999
$ws.Range("N83").Value = -27042.125
$ws.Range("H102").Value = 41669944
$ws.Range("I102").Value = 50003732
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 50003732
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = -50002110
$ws.Range("N102").Value = -4244
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344
$ws.Range("H111").Value = 20000
$ws.Range("J111").Value = 20000
$ws.Range("L111").Value = 20000
$ws.Range("N111").Value = -26134
$ws.Range("H122").Value = 41667650
$ws.Range("I122").Value = 15874097
$ws.Range("J122").Value = 90909880
$ws.Range("K122").Value = 47622291
$ws.Range("L122").Value = 272729640
$ws.Range("M122").Value = -47619841
$ws.Range("N122").Value = -272734540
$ws.Range("H132").Value = 19051.5
$ws.Range("I132").Value = 3439.8
$ws.Range("K132").Value = 10319.4
$ws.Range("M132").Value = -7789.400000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2752.9092
$ws.Range("I7").Value = 3501.2307
$ws.Range("J7").Value = 2266.5
$ws.Range("K7").Value = 3501.2307
$ws.Range("L7").Value = 2266.5
$ws.Range("M7").Value = -3389.2307
$ws.Range("N7").Value = -2490.5
$ws.Range("H22").Value = 1997.15
$ws.Range("I22").Value = 1835.8
$ws.Range("J22").Value = 2481.2
$ws.Range("K22").Value = 1835.8
$ws.Range("L22").Value = 2481.2
$ws.Range("M22").Value = -1540.8
$ws.Range("N22").Value = -3071.2
$ws.Range("H27").Value = 1997.15
$ws.Range("I27").Value = 1835.8
$ws.Range("J27").Value = 2481.2
$ws.Range("K27").Value = 1835.8
$ws.Range("L27").Value = 2481.2
$ws.Range("M27").Value = -1728.8
$ws.Range("N27").Value = -2695.2
$ws.Range("H46").Value = 748.5714
$ws.Range("I46").Value = 636.375
$ws.Range("J46").Value = 898.1667
$ws.Range("K46").Value = 636.375
$ws.Range("L46").Value = 898.1667
$ws.Range("M46").Value = -448.375
$ws.Range("N46").Value = -1274.1667
$ws.Range("H110").Value = 37999.5
$ws.Range("J110").Value = 37999.5
$ws.Range("L110").Value = 37999.5
$ws.Range("N110").Value = -46179.5
$ws.Range("H126").Value = 2752.9092
$ws.Range("I126").Value = 3501.2307
$ws.Range("J126").Value = 2266.5
$ws.Range("K126").Value = 10503.6921
$ws.Range("L126").Value = 6799.5
$ws.Range("M126").Value = -8033.6921
$ws.Range("N126").Value = -11739.5
$ws.Range("H132").Value = 2378.2163
$ws.Range("I132").Value = 1113.4546
$ws.Range("J132").Value = 4233.2
$ws.Range("K132").Value = 3340.3638
$ws.Range("L132").Value = 12699.6
$ws.Range("M132").Value = -810.3638000000001
$ws.Range("N132").Value = -17759.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 26184
$ws.Range("J92").Value = 26184
$ws.Range("L92").Value = 26184
$ws.Range("N92").Value = -31176
$ws.Range("H132").Value = 1730.125
$ws.Range("I132").Value = 1244
$ws.Range("J132").Value = 2799.6
$ws.Range("K132").Value = 3732
$ws.Range("L132").Value = 8398.799999999999
$ws.Range("M132").Value = -1202
$ws.Range("N132").Value = -13458.8
$ws.Range("H136").Value = 30304824
$ws.Range("I136").Value = 45456052
$ws.Range("J136").Value = 2373.0908
$ws.Range("K136").Value = 136368156
$ws.Range("L136").Value = 7119.2724
$ws.Range("M136").Value = -136365606
$ws.Range("N136").Value = -12219.2724
